$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the confidential disclosure text date (2021-05-20 -> 2021-05-21)
$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-05-21 for illustrative purposes only and are subject to change."

# Update the Weight / Percent Change values for rows 2-6
$ws.Range("D2").Value = 0.2555663594986085
$ws.Range("E2").Value = 0.002198819581067157

$ws.Range("D3").Value = 0.2539884317022559
$ws.Range("E3").Value = 0.009943563558183044

$ws.Range("D4").Value = 0.2431367223363745
$ws.Range("E4").Value = 0.004876145894284978

$ws.Range("D5").Value = 0.2473084864627612
$ws.Range("E5").Value = -0.003371369294605908

$ws.Range("E6").Value = 0.003439296322016672
